$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of mod-count tracking data for 2025/12/02, appended after row 22.
# Match the center/center alignment style used by the existing data rows.
$ws.Range("A23:C23").HorizontalAlignment = -4108
$ws.Range("A23:C23").VerticalAlignment = -4108

# Force the date-like text to be stored as literal text (not auto-converted
# to a date serial number) by formatting the cell as Text before entry, then
# reverting the display format back to General afterwards.
$ws.Range("A23").NumberFormat = "@"
$ws.Range("A23").Value = "2025/12/02"
$ws.Range("A23").NumberFormat = "general"

$ws.Range("B23").Value = "逃离鸭科夫"
$ws.Range("C23").Value = 1310
